$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "All.deja.sex"
$ws.Range("C1").Value = "Males.deja.sex"
$ws.Range("D1").Value = "Females.deja.sex"
$ws.Range("E1").Value = "Not known / missing.deja.sex"
